# Refresh cryptocurrency price / 1h-volume figures on Sheet1 (cryptos.xlsx),
# mirroring the scheduled "Updated cryptos list" GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.888.31'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.617.60'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').Value = "'212.41"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('D6').Value = "'0.499"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('E9').Value = '  -1.11%  '
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.841.72'
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('D13').Value = '1.605.72'
$ws.Range('E13').Value = '  -9.17%  '
$ws.Range('D14').Value = "'4.12"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').Value = '25.897.21'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '0.0₃0736'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').Value = "'190.90"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = "'4.24"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('E23').Value = '  -1.90%  '
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('D27').Value = "'1.72"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.92%  '
$ws.Range('D28').Value = "'6.63"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.72%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = "'1.23"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  -1.70%  '
$ws.Range('D33').Value = "'3.09"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('D36').Value = '1.125.02'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  -4.59%  '
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = "'0.510"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0153"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.31%  '
$ws.Range('D41').Value = "'97.81"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').Value = '1.752.68'
$ws.Range('D43').Value = "'0.748"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.73%  '
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('D47').Value = "'53.98"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.14%  '
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('E51').Value = '  -0.49%  '
